# Scheduled price/profit refresh for the Fenrir_Profits leve tables.
# Updates market-price-derived columns (H:N) on affected rows across
# the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 5545.1577
$ws.Cells.Item(69, 9).Value = 4504.3335
$ws.Cells.Item(69, 10).Value = 5740.3125
$ws.Cells.Item(69, 11).Value = 13513.0005
$ws.Cells.Item(69, 12).Value = 17220.9375
$ws.Cells.Item(69, 13).Value = -12639.0005
$ws.Cells.Item(69, 14).Value = -18968.9375
$ws.Cells.Item(72, 8).Value = 5545.1577
$ws.Cells.Item(72, 9).Value = 4504.3335
$ws.Cells.Item(72, 10).Value = 5740.3125
$ws.Cells.Item(72, 11).Value = 40539.0015
$ws.Cells.Item(72, 12).Value = 51662.8125
$ws.Cells.Item(72, 13).Value = -36171.0015
$ws.Cells.Item(72, 14).Value = -60398.8125
$ws.Cells.Item(96, 8).Value = 529.1111
$ws.Cells.Item(96, 9).Value = 324
$ws.Cells.Item(96, 10).Value = 734.2222
$ws.Cells.Item(96, 11).Value = 972
$ws.Cells.Item(96, 12).Value = 2202.6666
$ws.Cells.Item(96, 13).Value = 401
$ws.Cells.Item(96, 14).Value = -4948.6666
$ws.Cells.Item(127, 8).Value = 2279.75
$ws.Cells.Item(127, 9).Value = 1004.4286
$ws.Cells.Item(127, 10).Value = 2966.4614
$ws.Cells.Item(127, 11).Value = 3013.2858
$ws.Cells.Item(127, 12).Value = 8899.3842
$ws.Cells.Item(127, 13).Value = 1946.7142
$ws.Cells.Item(127, 14).Value = -18819.3842
$ws.Cells.Item(129, 8).Value = 714.56525
$ws.Cells.Item(129, 10).Value = 1036.0834
$ws.Cells.Item(129, 12).Value = 3108.2502
$ws.Cells.Item(129, 14).Value = -13108.2502
$ws.Cells.Item(132, 8).Value = 68255490
$ws.Cells.Item(132, 9).Value = 75440010
$ws.Cells.Item(132, 10).Value = 2550
$ws.Cells.Item(132, 11).Value = 226320030
$ws.Cells.Item(132, 12).Value = 7650
$ws.Cells.Item(132, 13).Value = -226317500
$ws.Cells.Item(132, 14).Value = -12710
$ws.Cells.Item(135, 8).Value = 6685.7144
$ws.Cells.Item(135, 9).Value = 9446.77
$ws.Cells.Item(135, 10).Value = 2199
$ws.Cells.Item(135, 11).Value = 85020.93000000001
$ws.Cells.Item(135, 12).Value = 19791
$ws.Cells.Item(135, 13).Value = -82485.93000000001
$ws.Cells.Item(135, 14).Value = -24861
$ws.Cells.Item(137, 8).Value = 592004.1
$ws.Cells.Item(137, 9).Value = 916812.75
$ws.Cells.Item(137, 10).Value = 64190
$ws.Cells.Item(137, 11).Value = 2750438.25
$ws.Cells.Item(137, 12).Value = 192570
$ws.Cells.Item(137, 13).Value = -2747888.25
$ws.Cells.Item(137, 14).Value = -197670
$ws.Cells.Item(138, 8).Value = 1720.3654
$ws.Cells.Item(138, 9).Value = 1053.1923
$ws.Cells.Item(138, 10).Value = 2387.5386
$ws.Cells.Item(138, 11).Value = 3159.5769
$ws.Cells.Item(138, 12).Value = 7162.6158
$ws.Cells.Item(138, 13).Value = 1980.4231
$ws.Cells.Item(138, 14).Value = -17442.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 10264.25
$ws.Cells.Item(61, 9).Value = 11159.143
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 11159.143
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 13).Value = -10947.143
$ws.Cells.Item(61, 14).Value = -4424
$ws.Cells.Item(132, 8).Value = 3126523.8
$ws.Cells.Item(132, 9).Value = 4033246.8
$ws.Cells.Item(132, 10).Value = 3366.3333
$ws.Cells.Item(132, 11).Value = 12099740.4
$ws.Cells.Item(132, 12).Value = 10098.9999
$ws.Cells.Item(132, 13).Value = -12097210.4
$ws.Cells.Item(132, 14).Value = -15158.9999
$ws.Cells.Item(136, 8).Value = 10264.25
$ws.Cells.Item(136, 9).Value = 11159.143
$ws.Cells.Item(136, 10).Value = 4000
$ws.Cells.Item(136, 11).Value = 33477.429
$ws.Cells.Item(136, 12).Value = 12000
$ws.Cells.Item(136, 13).Value = -30927.429
$ws.Cells.Item(136, 14).Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 11513361
$ws.Cells.Item(134, 9).Value = 17571850
$ws.Cells.Item(134, 10).Value = 2231.8
$ws.Cells.Item(134, 11).Value = 52715550
$ws.Cells.Item(134, 12).Value = 6695.400000000001
$ws.Cells.Item(134, 13).Value = -52713015
$ws.Cells.Item(134, 14).Value = -11765.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 14851.2705
$ws.Cells.Item(31, 9).Value = 5086.7915
$ws.Cells.Item(31, 10).Value = 32878
$ws.Cells.Item(31, 11).Value = 5086.7915
$ws.Cells.Item(31, 12).Value = 32878
$ws.Cells.Item(31, 13).Value = -4791.7915
$ws.Cells.Item(31, 14).Value = -33468
$ws.Cells.Item(34, 8).Value = 14851.2705
$ws.Cells.Item(34, 9).Value = 5086.7915
$ws.Cells.Item(34, 10).Value = 32878
$ws.Cells.Item(34, 11).Value = 5086.7915
$ws.Cells.Item(34, 12).Value = 32878
$ws.Cells.Item(34, 13).Value = -4884.7915
$ws.Cells.Item(34, 14).Value = -33282
$ws.Cells.Item(58, 8).Value = 7159036
$ws.Cells.Item(58, 9).Value = 10205300
$ws.Cells.Item(58, 10).Value = 51088
$ws.Cells.Item(58, 11).Value = 10205300
$ws.Cells.Item(58, 12).Value = 51088
$ws.Cells.Item(58, 13).Value = -10205097
$ws.Cells.Item(58, 14).Value = -51494
$ws.Cells.Item(99, 8).Value = 38462524
$ws.Cells.Item(99, 9).Value = 58824388
$ws.Cells.Item(99, 10).Value = 1218.2222
$ws.Cells.Item(99, 11).Value = 58824388
$ws.Cells.Item(99, 12).Value = 1218.2222
$ws.Cells.Item(99, 13).Value = -58822890
$ws.Cells.Item(99, 14).Value = -4214.2222
$ws.Cells.Item(126, 8).Value = 38462524
$ws.Cells.Item(126, 9).Value = 58824388
$ws.Cells.Item(126, 10).Value = 1218.2222
$ws.Cells.Item(126, 11).Value = 176473164
$ws.Cells.Item(126, 12).Value = 3654.6666
$ws.Cells.Item(126, 13).Value = -176470694
$ws.Cells.Item(126, 14).Value = -8594.6666
$ws.Cells.Item(132, 8).Value = 11499103
$ws.Cells.Item(132, 9).Value = 18519396
$ws.Cells.Item(132, 10).Value = 11351.091
$ws.Cells.Item(132, 11).Value = 55558188
$ws.Cells.Item(132, 12).Value = 34053.273
$ws.Cells.Item(132, 13).Value = -55555658
$ws.Cells.Item(132, 14).Value = -39113.273
$ws.Cells.Item(134, 8).Value = 11576391
$ws.Cells.Item(134, 9).Value = 20836248
$ws.Cells.Item(134, 10).Value = 4168504.5
$ws.Cells.Item(134, 11).Value = 62508744
$ws.Cells.Item(134, 12).Value = 12505513.5
$ws.Cells.Item(134, 13).Value = -62506209
$ws.Cells.Item(134, 14).Value = -12510583.5
$ws.Cells.Item(136, 8).Value = 7159036
$ws.Cells.Item(136, 9).Value = 10205300
$ws.Cells.Item(136, 10).Value = 51088
$ws.Cells.Item(136, 11).Value = 30615900
$ws.Cells.Item(136, 12).Value = 153264
$ws.Cells.Item(136, 13).Value = -30613350
$ws.Cells.Item(136, 14).Value = -158364

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 701.3333
$ws.Cells.Item(5, 9).Value = 701.3333
$ws.Cells.Item(5, 11).Value = 2103.9999
$ws.Cells.Item(5, 13).Value = -1991.9999
$ws.Cells.Item(22, 8).Value = 743.2558
$ws.Cells.Item(22, 10).Value = 743.2558
$ws.Cells.Item(22, 12).Value = 2229.7674
$ws.Cells.Item(22, 14).Value = -2567.7674
$ws.Cells.Item(27, 8).Value = 743.2558
$ws.Cells.Item(27, 10).Value = 743.2558
$ws.Cells.Item(27, 12).Value = 2229.7674
$ws.Cells.Item(27, 14).Value = -2433.7674
$ws.Cells.Item(135, 8).Value = 701.3333
$ws.Cells.Item(135, 9).Value = 701.3333
$ws.Cells.Item(135, 11).Value = 6311.9997
$ws.Cells.Item(135, 13).Value = -3776.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1245.6111
$ws.Cells.Item(102, 9).Value = 1182.5625
$ws.Cells.Item(102, 11).Value = 1182.5625
$ws.Cells.Item(102, 13).Value = 439.4375
$ws.Cells.Item(122, 8).Value = 1242.0834
$ws.Cells.Item(122, 9).Value = 1078.1111
$ws.Cells.Item(122, 10).Value = 1734
$ws.Cells.Item(122, 11).Value = 3234.3333
$ws.Cells.Item(122, 12).Value = 5202
$ws.Cells.Item(122, 13).Value = -784.3333000000002
$ws.Cells.Item(122, 14).Value = -10102
$ws.Cells.Item(132, 8).Value = 35719020
$ws.Cells.Item(132, 9).Value = 71429980
$ws.Cells.Item(132, 10).Value = 8051.7856
$ws.Cells.Item(132, 11).Value = 214289940
$ws.Cells.Item(132, 12).Value = 24155.3568
$ws.Cells.Item(132, 13).Value = -214287410
$ws.Cells.Item(132, 14).Value = -29215.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 9820126
$ws.Cells.Item(132, 9).Value = 22223502
$ws.Cells.Item(132, 10).Value = 2843227.5
$ws.Cells.Item(132, 11).Value = 66670506
$ws.Cells.Item(132, 12).Value = 8529682.5
$ws.Cells.Item(132, 13).Value = -66667976
$ws.Cells.Item(132, 14).Value = -8534742.5
